$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marks")

# --- Row 5: grades for Ruben Rubio Del Castillo ---
# Columns: B/C = Practica 0 (Tiempos), D/E = Practica 2, F/G = Practica 3,
#          H/I = Practica 4, J/K = Practica 5, L/M = Practica 6

# Copy the formatting already used by the "comment" cells (D5) onto the
# newly-graded comment cells (B5, J5, L5), and the formatting used by the
# "score" cells (E5) onto the newly-graded score cells (C5, K5, M5).
$ws.Range("D5").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("J5").PasteSpecial(-4122)
$ws.Range("L5").PasteSpecial(-4122)

$ws.Range("E5").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("K5").PasteSpecial(-4122)
$ws.Range("M5").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Practica 0 (Tiempos, columns B/C): comment + numeric mark
$ws.Range("B5").Value = "No incorpora los PDF. Tablas bien, pero no hay comentarios ni gráficas. Ni unidades de tiempo."
$ws.Range("C5").Value = 5

# Practica 5 (Prog. dinamica, columns J/K): comment + mark (asterisk = pending)
$ws.Range("J5").Value = "Método recursivo bien. PD al revés y añade un * y se come caracteres, la tabla de valores intermedios tiene las coordenadas a 0. Código no implementa maximo. Tablas de tiempos dentro de los valores permitidos, pero faltan datos pedidos: unidades de tiempos, procesador donde se mide y no hay comentarios."
$ws.Range("K5").Value = "*"

# Practica 6 (Backtracking, columns L/M): comment only, mark still pending (blank)
$ws.Range("L5").Value = "No tiene main en MejorLista para probar con distintos casos."
$ws.Range("L5").Font.Color = 255
$ws.Range("M5").Value = ""
